$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells text value while preserving its original cell
# style/format (Excel would otherwise auto-convert numeric-looking
# strings like "1.00" or "67.959.40" into actual numbers).
function Set-TextValue($addr, $value) {
  $range = $ws.Range($addr)
  $origStyle = $range.Style
  $range.NumberFormat = "@"
  $range.Value = $value
  $range.Style = $origStyle
}

Set-TextValue "D2" '67.959.40'
Set-TextValue "E2" '  +1.70%  '

Set-TextValue "D3" '3.330.36'
Set-TextValue "E3" '  +1.45%  '

Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.06%  '

Set-TextValue "D5" '581.92'
Set-TextValue "E5" '  +1.50%  '

Set-TextValue "D6" '177.61'
Set-TextValue "E6" '  +1.44%  '

Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  -0.03%  '

Set-TextValue "E8" '  +1.71%  '

Set-TextValue "D9" '3.327.32'
Set-TextValue "E9" '  +1.47%  '

Set-TextValue "D10" '0.184'
Set-TextValue "E10" '  +5.98%  '

Set-TextValue "D11" '0.582'
Set-TextValue "E11" '  +1.60%  '

Set-TextValue "D12" '47.10'
Set-TextValue "E12" '  +3.72%  '

Set-TextValue "E13" '  +2.34%  '

Set-TextValue "D14" '682.95'
Set-TextValue "E14" '  -0.23%  '

Set-TextValue "D15" '3.872.34'
Set-TextValue "E15" '  +1.57%  '

Set-TextValue "E16" '  +1.74%  '

Set-TextValue "D17" '67.949.17'
Set-TextValue "E17" '  +1.44%  '

Set-TextValue "B18" 'WrappedEther'
Set-TextValue "C18" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D18" '3.346.42'
Set-TextValue "E18" '  +1.93%  '

Set-TextValue "B19" 'TRON'
Set-TextValue "C19" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D19" '0.118'
Set-TextValue "E19" '  -0.44%  '

Set-TextValue "D20" '17.45'
Set-TextValue "E20" '  +1.07%  '

Set-TextValue "D21" '11.07'
Set-TextValue "E21" '  +3.39%  '

Set-TextValue "E22" '  +1.17%  '

Set-TextValue "D23" '5.40'
Set-TextValue "E23" '  +5.01%  '

Set-TextValue "D24" '17.05'
Set-TextValue "E24" '  +0.43%  '

Set-TextValue "D25" '99.66'
Set-TextValue "E25" '  +1.08%  '

Set-TextValue "E26" '  +1.13%  '

Set-TextValue "E27" '  +0.08%  '

Set-TextValue "D28" '9.58'
Set-TextValue "E28" '  +3.51%  '

Set-TextValue "D29" '33.19'
Set-TextValue "E29" '  -1.03%  '

Set-TextValue "E30" '  +2.55%  '

Set-TextValue "D31" '7.11'
Set-TextValue "E31" '  +5.56%  '

Set-TextValue "D32" '565.91'
Set-TextValue "E32" '  -0.74%  '

Set-TextValue "D33" '11.01'
Set-TextValue "E33" '  +1.76%  '

Set-TextValue "E34" '  +2.71%  '

Set-TextValue "B35" 'Dai'
Set-TextValue "C35" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D35" '1.00'
Set-TextValue "E35" '  -0.09%  '

Set-TextValue "B36" 'OKB'
Set-TextValue "C36" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D36" '57.17'
Set-TextValue "E36" '  +3.49%  '

Set-TextValue "D37" '3.705.62'
Set-TextValue "E37" '  -4.35%  '

Set-TextValue "D38" '3.39'
Set-TextValue "E38" '  +2.86%  '

Set-TextValue "E39" '  +8.83%  '

Set-TextValue "D40" '0.132'
Set-TextValue "E40" '  +3.89%  '

Set-TextValue "E41" '  +6.90%  '

Set-TextValue "E42" '  +2.65%  '

Set-TextValue "B43" 'PEPE'
Set-TextValue "C43" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D43" '0.0₃0675'
Set-TextValue "E43" '  +1.39%  '

Set-TextValue "B44" 'ApeXProtocol'
Set-TextValue "C44" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D44" '3.34'
Set-TextValue "E44" '  -1.16%  '

Set-TextValue "D45" '0.337'
Set-TextValue "E45" '  +3.51%  '

Set-TextValue "E46" '  +0.74%  '

Set-TextValue "E47" '  +4.88%  '

Set-TextValue "E48" '  +1.31%  '

Set-TextValue "E49" '  -0.38%  '

Set-TextValue "E50" '  -2.81%  '

Set-TextValue "D51" '130.30'
